# Update computed results on the "Gas Sorption Input" sheet (P12:Q18)
# with values produced by a more efficient loading routine.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gas Sorption Input")

$ws.Range("P12").Value = 0.0375898586132879
$ws.Range("Q12").Value = 0.0005889352803773407

$ws.Range("P13").Value = 0.05867943939644259
$ws.Range("Q13").Value = 0.0005068691481518938

$ws.Range("P14").Value = 0.06968291800307592
$ws.Range("Q14").Value = 0.0008593573285447139

$ws.Range("P15").Value = 0.07752732911814221
$ws.Range("Q15").Value = 0.0017373645222090906

$ws.Range("P16").Value = 0.08424243429322628
$ws.Range("Q16").Value = 0.0027979286560536696

$ws.Range("P17").Value = 0.09068008289072088
$ws.Range("Q17").Value = 0.004013519039213339

$ws.Range("P18").Value = 0.09613440530926486
$ws.Range("Q18").Value = 0.005198235451499367
